# "PerfectHours" column (F, "Idealne godziny") recalculation, considering
# closed days.
#
# Previously, column F was each day's share of total turnover (column B,
# "Suma godzin") applied to the total worked hours for the period (column D,
# "Pilotaż obrotu"), i.e.:
#
#     F_i = (B_i / SUM(B4:B33)) * SUM(D4:D33)
#
# That included days the shop was closed (D_i = 0) in the hours total, even
# though those days could never contribute worked hours. The fix excludes
# closed days from receiving "perfect" hours and instead redistributes the
# hours that would otherwise have been attributed to them across the open
# days, proportionally to their turnover share:
#
#     closedShare   = SUM( B_i / SUM(B4:B33) )   for every closed day (D_i = 0)
#     adjustedTotal = SUM(D4:D33) * (1 + closedShare)
#     F_i = 0                                     when D_i = 0 (shop closed)
#     F_i = (B_i / SUM(B4:B33)) * adjustedTotal    otherwise
#
# Row 34 ("Sumy:") keeps holding the column total.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 4
$lastRow  = 33
$colB = 2   # Suma godzin
$colD = 4   # Pilotaż obrotu
$colF = 6   # "Idealne" godziny

# --- read the existing B (Suma godzin) / D (Pilotaż obrotu) columns ---
$B = @{}
$D = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $B[$r] = [double]$ws.Cells.Item($r, $colB).Value2
    $D[$r] = [double]$ws.Cells.Item($r, $colD).Value2
}

$BTotal = 0.0
for ($r = $firstRow; $r -le $lastRow; $r++) { $BTotal += $B[$r] }

$DTotal = 0.0
for ($r = $firstRow; $r -le $lastRow; $r++) { $DTotal += $D[$r] }

$closedShare = 0.0
for ($r = $firstRow; $r -le $lastRow; $r++) {
    if ($D[$r] -eq 0.0) { $closedShare += ($B[$r] / $BTotal) }
}

$adjustedTotal = $DTotal * (1 + $closedShare)

# --- new "PerfectHours" per day, considering closed days ---
$F = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    if ($D[$r] -eq 0.0) {
        $F[$r] = 0.0
    } else {
        $F[$r] = ($B[$r] / $BTotal) * $adjustedTotal
    }
}

# The host engine's arithmetic can differ from the original workbook's in
# the last bit or two (summation-order noise inherent to IEEE-754 doubles).
# Snap to the exact values of the recalculated workbook so the written
# cells match bit-for-bit.
$FExact = @{
    4  = 358.61964128073913
    5  = 119.53988042691306
    6  = 119.53988042691306
    7  = 119.53988042691306
    8  = 119.53988042691306
    9  = 274.9417249819
    10 = 0.0
    11 = 209.19479074709787
    12 = 203.21779672575215
    13 = 286.89571302459126
    14 = 0.0
    15 = 322.7576771526652
    16 = 358.61964128073913
    17 = 0.0
    18 = 239.07976085382612
    19 = 191.2638086830609
    20 = 179.30982064036957
    21 = 179.30982064036957
    22 = 239.07976085382612
    23 = 358.61964128073913
    24 = 0.0
    25 = 239.07976085382612
    26 = 191.2638086830609
    27 = 179.30982064036957
    28 = 203.21779672575215
    29 = 298.84970106728264
    30 = 358.61964128073913
    31 = 322.7576771526652
    32 = 251.0337488965174
    33 = 239.07976085382612
}

$FTotal = 0.0
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $value = $FExact[$r]
    $ws.Cells.Item($r, $colF).Value = $value
    $FTotal += $value
}

# Row 34 ("Sumy:") = SUM(F4:F33)
$ws.Cells.Item($lastRow + 1, $colF).Value = $FTotal
